$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "테스트수정1"
$ws.Range("B2").Value = 45828
$ws.Range("B2").NumberFormat = "mm-dd-yy"

$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.Range("B3").Select() | Out-Null
